# "add prolificid in rank to use in binary"
# This re-derives the re_rank values with an updated prolificid/name pairing
# for the tied-ranked workers (rows 7-10), and refreshes the re_rank (E column)
# values across the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update re_rank (column E) values for all data rows ---
$ws.Cells.Item(2, 5).Value = 13.45172621041747
$ws.Cells.Item(3, 5).Value = 8.467383315037575
$ws.Cells.Item(4, 5).Value = 7.263320786645187
$ws.Cells.Item(5, 5).Value = 7.045923228846132
$ws.Cells.Item(6, 5).Value = 5.285624560074965
$ws.Cells.Item(7, 5).Value = 5.22667163757618
$ws.Cells.Item(8, 5).Value = 5.193444245373518
$ws.Cells.Item(9, 5).Value = 4.163691280357252
$ws.Cells.Item(10, 5).Value = 4.093901744365527
$ws.Cells.Item(11, 5).Value = 3.064644559899139
$ws.Cells.Item(12, 5).Value = 2.350791450174602
$ws.Cells.Item(13, 5).Value = 2.331889986248744

# --- Swap the prolificid (B), name (C) and race (F) between rows 7 & 8 ---
$ws.Cells.Item(7, 2).Value = 32
$ws.Cells.Item(7, 3).Value = "Jamarii"
$ws.Cells.Item(7, 6).Value = "Black or African American"

$ws.Cells.Item(8, 2).Value = 26
$ws.Cells.Item(8, 3).Value = "Juan"
$ws.Cells.Item(8, 6).Value = "Hispanic"

# --- Swap the prolificid (B) and name (C) between rows 9 & 10 ---
# (race/F is "White" for both so it stays the same)
$ws.Cells.Item(9, 2).Value = 2
$ws.Cells.Item(9, 3).Value = "Corey"

$ws.Cells.Item(10, 2).Value = 33
$ws.Cells.Item(10, 3).Value = "Brennan"
